$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -----------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- Shared style for B1 / A2: bold font, thin box border, --------------
# --- centered horizontally, top-aligned vertically ----------------------
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108   # xlCenter
$b1.VerticalAlignment = -4160     # xlTop
$b1.Borders.LineStyle = 1         # xlContinuous
$b1.Borders.Weight = 2            # xlThin

# Apply the exact same style to A2 by copying B1's format, so both cells
# share a single, de-duplicated cell-format entry.
$b1.Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
